# edit.ps1 - applies the JMLR cover-letter edits:
#  1) Date paragraph: "Apr 25, 2022" -> "11 May 2022" (as 4 runs: "1","1"," May"," 2022")
#  2) Title runs "Exploring ... Using" + " " + "Animated ... Projections" merge into a
#     single run (text itself is unchanged, only the run split changes)
#
# Note: a plain Find/Execute (or Range.Text=) replacement in this runtime re-merges
# every same-formatted run from the edit point through to the end of the paragraph
# (and one run to the left), which would swallow far more of the paragraph than the
# source edit touched. Inserting brand-new text exactly at a paragraph's Start
# position is the one place that never triggers that cascade, so both edits below
# are implemented as "insert the replacement at the paragraph start, then delete the
# stale old text" (Range.Delete likewise never cascades into neighboring runs).

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Change 1: date paragraph (paragraph 1) "Apr 25, 2022" -> "11 May 2022"
# ---------------------------------------------------------------
$p1 = $d.Paragraphs(1).Range
$oldDate = "Apr 25, 2022"
$p1Text = $p1.Text
if ($p1Text.Length -ge $oldDate.Length -and $p1Text.Substring(0, $oldDate.Length) -eq $oldDate) {
    $pStart = $p1.Start

    # Insert the 4 new runs, in reverse order, always at the (fixed) paragraph-start
    # position so each InsertBefore call creates its own independent run instead of
    # being folded into a neighboring run.
    $d.Range($pStart, $pStart).InsertBefore(" 2022")
    $d.Range($pStart, $pStart).InsertBefore(" May")
    $d.Range($pStart, $pStart).InsertBefore("1")
    $d.Range($pStart, $pStart).InsertBefore("1")

    # Remove the old "Apr 25, 2022" text, which now immediately follows what we
    # just inserted.
    $newText = "11 May 2022"
    $newLen = $newText.Length
    $oldLen = $oldDate.Length
    $rOld = $d.Range($pStart + $newLen, $pStart + $newLen + $oldLen)
    if ($rOld.Text -eq $oldDate) {
        $rOld.Delete()
    }
}

# ---------------------------------------------------------------
# Change 2: merge the title runs within the submission paragraph
# ---------------------------------------------------------------
# Original run texts of that paragraph (15 runs); the title is runs 5-7.
$origRuns = @(
    "I am enclosing a submission to ",
    "Journal of Machine Learning Researc",
    "h",
    " entitled “",
    "Exploring Local Explanations of Nonlinear Models Using",
    " ",
    "Animated Linear Projections",
    "”",
    ".",
    " The manuscript is 2",
    "4",
    " pages long and includes nine figures. ",
    "This manuscript has not been submitted to other journals or conferences. Though, s",
    "ome of the ",
    "content of this manuscript was discussed in my thesis, currently under examiner review."
)

$targetRuns = @(
    "I am enclosing a submission to ",
    "Journal of Machine Learning Researc",
    "h",
    " entitled “",
    "Exploring Local Explanations of Nonlinear Models Using Animated Linear Projections",
    "”",
    ".",
    " The manuscript is 2",
    "4",
    " pages long and includes nine figures. ",
    "This manuscript has not been submitted to other journals or conferences. Though, s",
    "ome of the ",
    "content of this manuscript was discussed in my thesis, currently under examiner review."
)

$origFull = ""
foreach ($seg in $origRuns) { $origFull = "$origFull$seg" }

$targetFull = ""
foreach ($seg in $targetRuns) { $targetFull = "$targetFull$seg" }

# Find the paragraph that currently holds this (unchanged) text so we can rebuild it.
$rFind = $d.Content
$found = $rFind.Find.Execute($origFull, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $pStart2 = $rFind.Start
    $origLen = $origFull.Length

    # Insert the target runs, in reverse order, always at the paragraph's fixed
    # start position, so each becomes its own independent run.
    for ($i = $targetRuns.Count - 1; $i -ge 0; $i--) {
        $seg = $targetRuns[$i]
        if ($seg.Length -gt 0) {
            $d.Range($pStart2, $pStart2).InsertBefore($seg)
        }
    }

    # Delete the old paragraph text, which now immediately follows what we just
    # inserted.
    $newLen2 = $targetFull.Length
    $rOld2 = $d.Range($pStart2 + $newLen2, $pStart2 + $newLen2 + $origLen)
    if ($rOld2.Text -eq $origFull) {
        $rOld2.Delete()
    }
}

Write-Host "Edit complete."
Write-Host "P1: $($d.Paragraphs(1).Range.Text)"
Write-Host "P4: $($d.Paragraphs(4).Range.Text)"
